$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Conta 004479965 / DIEGO / 17432.65  ->  004482102 / NATALIA / 37567.8
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "004482102"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "NATALIA"
$ws.Range("C3").Value = 37567.8

# Row 4: Conta 005428871 / ROSANGELA / 14879.47 -> 005002457 / ROSANGELA / 34484.08
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "005002457"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "ROSANGELA"
$ws.Range("C4").Value = 34484.08

# Row 5 (004267976 / E3 / 11892.73) is removed entirely, shifting all following rows up by one
$ws.Range("A5:C5").EntireRow.Delete()
